# Revert "adding term 2.0.0"
# - Metadata sheet: Version / Date / Contact values restored to the
#   pre-2.0.0 revision.
# - "Include from SNOMED CT" sheet: concept code restored to its
#   pre-2.0.0 value (160245001) and a second concept row (116223007) is
#   (re)inserted.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Metadata")
$ws2 = $wb.Worksheets.Item("Include from SNOMED CT")

# ---- Metadata sheet --------------------------------------------------
$ws1.Range("B3").Value = "1.1.0"
$ws1.Range("B8").Value = "2023-07-10T23:08:03+02:00"
$ws1.Range("B10").Value = "No display for ContactDetail"

# ---- "Include from SNOMED CT" sheet -----------------------------------
# Insert a new row 3 (pushes the old row3/row4 down to row4/row5).
$ws2.Rows("3:3").Insert()

# New row 3 should carry the same formatting as the surrounding data rows
# (row 2's style), not the blank default style Excel assigns on insert.
$ws2.Range("A2:B2").Copy()
$ws2.Range("A3:B3").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# Row 2 concept code: 116224001 -> 160245001
# Typing a purely numeric string directly would store it as a Number, so
# write it as a text formula then immediately collapse the formula to a
# plain cached value via a self copy/paste-values (keeps the cell's
# existing style and avoids leaving a live formula behind).
$ws2.Range("A2").Formula = '="160245001"'
$ws2.Range("A2").Copy()
$ws2.Range("A2").PasteSpecial(-4163)      # xlPasteValues
$excel.CutCopyMode = $false

# New row 3 concept code: 116223007
$ws2.Range("A3").Formula = '="116223007"'
$ws2.Range("A3").Copy()
$ws2.Range("A3").PasteSpecial(-4163)      # xlPasteValues
$excel.CutCopyMode = $false
$ws2.Range("B3").Value = ""
